# Sprint24Plan.xlsx - Build 452
# "Removed localized fields in account group"
#
# The sprint backlog on the "Backlog" sheet gets 4 new rows inserted right
# before the existing "آماده کردن اطلاعات گزارش مدیریت اسناد..." group
# (currently starting at row 10), describing the new "Account group" work
# item. Everything below shifts down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# 1) Insert 4 blank rows before row 10 - this pushes the existing rows
#    10-25 down to 14-29 and Excel auto-fixes dimension / data validation
#    ranges for us.
$ws.Rows("10:13").Insert()

# 2) The newly inserted rows don't inherit the table's row styling, so
#    copy the formatting (borders/alignment) from the row directly above
#    (row 9, which still has the regular data-row style) onto them.
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E13").PasteSpecial(-4122)
$ws.Range("A10:E13").RowHeight = 20.4

# 3) Fill in the new "Account group" backlog entries.
$ws.Range("A10").Value = "مدیریت اطلاعات گروه های حساب"
$ws.Range("B10").Value = "پیاده سازی عملیات اصلی مدیریت گروه های حساب در سرویس وب"
$ws.Range("C10").Value = 1

$ws.Range("B11").Value = "مدیریت اطلاعات گروه های حساب در برنامه وب"
$ws.Range("C11").Value = 1

$ws.Range("B12").Value = "اضافه کردن شناسه گروه حساب برای حساب های کل در سرویس وب"
$ws.Range("C12").Value = 1

$ws.Range("B13").Value = "امکان انتخاب گروه حساب هنگام ایجاد و اصلاح حساب های کل در برنامه وب"
$ws.Range("C13").Value = 1

# 4) Grow the existing table (Table1) so the new rows participate in it,
#    and its autofilter range grows with it.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:E29"))

# 5) Match the author's final selection.
$ws.Range("B13").Select()
